# Applies the Alvearie FHIR IG "wh-payer-coverage-class" ValueSet metadata
# update: version bump, refreshed date, Publisher value filled in, and the
# duplicated "Contact" row replaced by a "Jurisdiction" row (on the
# Metadata sheet only; the Codes sheets are left semantically unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: refreshed publish timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was a "Contact" row with no display text; it becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was the second, now-redundant "Contact" row -- remove it entirely,
# which shifts Description/Purpose/Copyright/Immutable up by one row.
$ws.Rows.Item(11).Delete()
